$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 370 (old rows 370-374 shift down to 373-377)
$ws.Range("A370:A372").EntireRow.Insert()

# New row 370: Especial / Brasil, week of 2021-09-09
$ws.Range("A370").Value = 8
$ws.Range("B370").Value = "Terminal La Palmera de La Serena"
$ws.Range("C370").Value = "Coquimbo"
$ws.Range("D370").Value = 44448
$ws.Range("E370").Value = 4
$ws.Range("F370").Value = "Fruta"
$ws.Range("G370").Value = 100108
$ws.Range("H370").Value = "Tropicales y subtropicales"
$ws.Range("I370").Value = 100108002
$ws.Range("J370").Value = "Mango"
$ws.Range("K370").Value = "Sin especificar"
$ws.Range("L370").Value = "Especial"
$ws.Range("M370").Value = 512
$ws.Range("N370").Value = 8000
$ws.Range("O370").Value = 8500
$ws.Range("P370").Value = 8250
$ws.Range("Q370").Value = "`$/bandeja 4 kilos"
$ws.Range("R370").Value = "Brasil"
$ws.Range("S370").Value = 2062
$ws.Range("T370").Value = 4

# New row 371: Primera / Brasil, week of 2021-09-09
$ws.Range("A371").Value = 8
$ws.Range("B371").Value = "Terminal La Palmera de La Serena"
$ws.Range("C371").Value = "Coquimbo"
$ws.Range("D371").Value = 44448
$ws.Range("E371").Value = 4
$ws.Range("F371").Value = "Fruta"
$ws.Range("G371").Value = 100108
$ws.Range("H371").Value = "Tropicales y subtropicales"
$ws.Range("I371").Value = 100108002
$ws.Range("J371").Value = "Mango"
$ws.Range("K371").Value = "Sin especificar"
$ws.Range("L371").Value = "Primera"
$ws.Range("M371").Value = 512
$ws.Range("N371").Value = 8000
$ws.Range("O371").Value = 8500
$ws.Range("P371").Value = 8250
$ws.Range("Q371").Value = "`$/bandeja 4 kilos"
$ws.Range("R371").Value = "Brasil"
$ws.Range("S371").Value = 2062
$ws.Range("T371").Value = 4

# New row 372: Segunda / Brasil, week of 2021-09-09
$ws.Range("A372").Value = 8
$ws.Range("B372").Value = "Terminal La Palmera de La Serena"
$ws.Range("C372").Value = "Coquimbo"
$ws.Range("D372").Value = 44448
$ws.Range("E372").Value = 4
$ws.Range("F372").Value = "Fruta"
$ws.Range("G372").Value = 100108
$ws.Range("H372").Value = "Tropicales y subtropicales"
$ws.Range("I372").Value = 100108002
$ws.Range("J372").Value = "Mango"
$ws.Range("K372").Value = "Sin especificar"
$ws.Range("L372").Value = "Segunda"
$ws.Range("M372").Value = 512
$ws.Range("N372").Value = 8000
$ws.Range("O372").Value = 8500
$ws.Range("P372").Value = 8250
$ws.Range("Q372").Value = "`$/bandeja 4 kilos"
$ws.Range("R372").Value = "Brasil"
$ws.Range("S372").Value = 2062
$ws.Range("T372").Value = 4
